# Updates on run.R for py files / Updates on time-series plots
#
# 01_admin/outline.xlsx gains a third column ("Script") recording which
# language each Verb/Folder step is implemented in. Every existing row is
# "R"; a new row is inserted after the "analyze / Counterfactual" entry for
# the equivalent "py" script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the "analyze / Counterfactual / py" entry -------
# Before the insert this is the 29th data row (sheet row 30): "analyze",
# "Counterfactual". We push it (and everything after it) down one row so a
# brand-new row 30 can hold the "py" counterpart.
$ws.Rows("30").Insert()

$ws.Range("A30").Value = "analyze"
$ws.Range("B30").Value = "Counterfactual"
$ws.Range("C30").Value = "py"

# --- New "Script" column ---------------------------------------------------
$ws.Range("C1").Value = "Script"

# Rows 2-29 kept their original data rows (1-28 of the table); fill "R".
$ws.Range("C2:C29").Value = "R"

# Rows 31-36 are the former rows 30-35, shifted down by the insert above.
$ws.Range("C31:C36").Value = "R"

# --- Selection / scroll position left by the editor ------------------------
$ws.Range("B30").Select()
